$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the candidate record on row 2 with refreshed test-credentials
$ws.Range("A2").Value = "pPKtD203"
$ws.Range("B2").Value = 23091447
$ws.Range("C2").Value = "rzljdex30"
$ws.Range("D2").Value = "Wt&5v2!B"
$ws.Range("F2").Value = "GznnOyiH"
$ws.Range("G2").Value = "wGTi"
